$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Original run text was "Modify the project". Split it into three runs:
#   "Modify the " / "project and " / "test build"
$tr.Text = "Modify the "
$tr.InsertAfter("project and ")
$tr.InsertAfter("test build")
